$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.5789666666666667
$ws.Range("H2").Value = 1.7369
$ws.Range("I2").Value = 0.01523705650035473
$ws.Range("J2").Value = 0.01523705650035472
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 44.63226435922223
$ws.Range("R2").Value = 401.690379233
$ws.Range("S2").Value = 0.003662709047407038
$ws.Range("T2").Value = 0.003662709047407038
$ws.Range("G3").Value = 0.5789666666666667
$ws.Range("H3").Value = 1.7369
$ws.Range("I3").Value = 0.01523705650035473
$ws.Range("J3").Value = 0.01523705650035472
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 58.81145561475557
$ws.Range("R3").Value = 529.3031005328
$ws.Range("S3").Value = 0.004826312392255611
$ws.Range("T3").Value = 0.004826312392255609
$ws.Range("G4").Value = 0.5789666666666667
$ws.Range("H4").Value = 1.7369
$ws.Range("I4").Value = 0.01523705650035473
$ws.Range("J4").Value = 0.01523705650035472
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 82.22877679768891
$ws.Range("R4").Value = 740.0589911792001
$ws.Range("S4").Value = 0.006748035060692078
$ws.Range("T4").Value = 0.006748035060692078
$ws.Range("I5").Value = 0.6545086962501954
$ws.Range("J5").Value = 0.6545086962501954
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 1917.181652228469
$ws.Range("R5").Value = 17254.63487005622
$ws.Range("S5").Value = 0.1573318917145425
$ws.Range("T5").Value = 0.1573318917145425
$ws.Range("I6").Value = 0.6545086962501954
$ws.Range("J6").Value = 0.6545086962501954
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.2073145447401761
$ws.Range("T6").Value = 0.2073145447401761
$ws.Range("I7").Value = 0.6545086962501954
$ws.Range("J7").Value = 0.6545086962501954
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.2898622597954768
$ws.Range("T7").Value = 0.2898622597954768
$ws.Range("I8").Value = 0.33025424724945
$ws.Range("J8").Value = 0.3302542472494499
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 967.3781066999545
$ws.Range("R8").Value = 8706.402960299591
$ws.Range("S8").Value = 0.0793870666107329
$ws.Range("T8").Value = 0.07938706661073289
$ws.Range("I9").Value = 0.33025424724945
$ws.Range("J9").Value = 0.3302542472494499
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.104607485445628
$ws.Range("T9").Value = 0.104607485445628
$ws.Range("I10").Value = 0.33025424724945
$ws.Range("J10").Value = 0.3302542472494499
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.1462596951930891
$ws.Range("T10").Value = 0.146259695193089
